# Fruta / hortaliza, semanal
# Inserts 6 new weekly price rows (2023-01-13, serial 44939) for
# "Vega Monumental Concepcion - Durazno" before the previous last
# entries (which get pushed down), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 259; this pushes the former
# rows 259-260 down to 265-266, preserving their original values.
$ws.Rows.Item(259).Resize(6).Insert()

# Common/static values shared by every new row in this block.
$mercadoId   = 11
$mercado     = "Vega Monumental Concepción"
$region      = "Bíobío"
$fecha       = 44939
$codreg      = 8
$tipo        = "Fruta"
$productoId  = 100103
$producto    = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria   = "Durazno"
$origen      = "Región de O'Higgins"

# Variable fields per new row (rows 259-264):
# Variedad, Calidad, N(envases), Min, Max, Prom, Unidad, PrecioKilo, Kilos
$newRows = @(
    @("Carson",      "Especial", 120, 15000, 15000, 15000, "$/caja 15 kilos empedrada", 1000, 15),
    @("Carson",      "Primera",  200, 14000, 14000, 14000, "$/caja 15 kilos empedrada",  933, 15),
    @("Carson",      "Segunda",  150, 13000, 13000, 13000, "$/caja 15 kilos empedrada",  867, 15),
    @("Springcrest",  "Especial", 150, 16000, 16000, 16000, "$/caja 15 kilos empedrada", 1067, 15),
    @("Springcrest",  "Primera",  200, 15000, 15000, 15000, "$/caja 15 kilos empedrada", 1000, 15),
    @("Springcrest",  "Segunda",  220, 13000, 13000, 13000, "$/caja 15 kilos empedrada",  867, 15)
)

$r = 259
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $row[0]
    $ws.Cells.Item($r, 12).Value = $row[1]
    $ws.Cells.Item($r, 13).Value = $row[2]
    $ws.Cells.Item($r, 14).Value = $row[3]
    $ws.Cells.Item($r, 15).Value = $row[4]
    $ws.Cells.Item($r, 16).Value = $row[5]
    $ws.Cells.Item($r, 17).Value = $row[6]
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $row[7]
    $ws.Cells.Item($r, 20).Value = $row[8]
    $r = $r + 1
}
